# Generate Report for Handoff
# Updates the localization status report: rows that were previously
# "Handed back: in sync with en-US" are now "Ready for handoff", and the
# handoff-generation timestamps are refreshed.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns + latest HO xliff generate date
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-09-03 07:03:16"

# zh-cn sheet: status + latest handoff datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-03 07:03:11"

# de-de sheet: status + latest handoff datetime
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-03 07:03:16"

# Columns auto-fit narrower now that "Ready for handoff" is shorter than
# "Handed back: in sync with en-US"
$overview.Range("E:F").ColumnWidth = 17.2159881591797
$zhcn.Range("C:C").ColumnWidth = 17.2159881591797
$dede.Range("C:C").ColumnWidth = 17.2159881591797
